# Show all upload errors: clear the two "missing data" cells on the
# "Direction 1 STOPS" sheet (B2/B3 -> blank, simulating rows that failed
# validation) and make that sheet the active one, matching the new
# validation flow that surfaces every row with missing data instead of
# stopping at the first one.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Direction 0 STOPS")
$ws2 = $wb.Worksheets.Item("Direction 1 STOPS")

# Remove the values that make this fixture "invalid" due to missing data.
$ws2.Range("B2").ClearContents()
$ws2.Range("B3").ClearContents()

# "Direction 1 STOPS" becomes the active/selected sheet & cell.
$ws2.Activate()
$ws2.Range("B3").Select()
